# Generate Report for Handoff
# Replace the old GUID-based file name with the new one across all sheets,
# and refresh the handoff/handback timestamps that accompany the new report.

$wb = $excel.ActiveWorkbook

$oldGuid = "f3b44bd0-28ce-4478-9fcd-c7a46703a62c"
$newGuid = "a6ead630-dd4b-499f-839c-c2dbc848ed83"

$oldZhXlf = "$oldGuid.5f1f4fd4a66eef46f5b0f7fe6ec5e3a5272aea77.zh-cn.xlf"
$newZhXlf = "$newGuid.b634111e2b80c3772bf5fb87b6ff1671ea2a5682.zh-cn.xlf"

$oldDeXlf = "$oldGuid.5f1f4fd4a66eef46f5b0f7fe6ec5e3a5272aea77.de-de.xlf"
$newDeXlf = "$newGuid.b634111e2b80c3772bf5fb87b6ff1671ea2a5682.de-de.xlf"

$oldHoDate = "2016-09-05 11:17:21"
$newHoDate = "2016-09-05 11:18:08"

$oldZhDate = "2016-09-05 11:17:06"
$newZhDate = "2016-09-05 11:17:53"

# ---- Overview sheet ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = "$newGuid.md"
$wsOverview.Range("B2").Value = "e2e\$newGuid.md"
$wsOverview.Range("G2").Value = $newHoDate

foreach ($hl in $wsOverview.Hyperlinks) {
    if ($hl.Range.Address() -eq '$B$2') {
        $hl.TextToDisplay = "e2e\$newGuid.md"
    }
}

# ---- zh-cn sheet ----
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("A2").Value = "$newGuid.md"
$wsZh.Range("G2").Value = $newZhXlf
$wsZh.Range("H2").Value = $newZhDate

foreach ($hl in $wsZh.Hyperlinks) {
    if ($hl.Range.Address() -eq '$A$2') {
        $hl.TextToDisplay = "$newGuid.md"
    }
}

# ---- de-de sheet ----
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("A2").Value = "$newGuid.md"
$wsDe.Range("G2").Value = $newDeXlf
$wsDe.Range("H2").Value = $newHoDate

foreach ($hl in $wsDe.Hyperlinks) {
    if ($hl.Range.Address() -eq '$A$2') {
        $hl.TextToDisplay = "$newGuid.md"
    }
}
